$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Update the valor mora (debt amount) value
$ws.Range("E11").Value = 252102

# Update period count
$ws.Range("F13").Value = 9

# Reorder the "Periodo Mora" values for the first worker (LUZ ESTELA SIERRA VELASQUEZ)
# from descending (1812..1805) to ascending (1805..1812)
$ws.Range("E16").Value = "1805"
$ws.Range("E17").Value = "1806"
$ws.Range("E18").Value = "1807"
$ws.Range("E19").Value = "1808"
$ws.Range("E20").Value = "1809"
$ws.Range("E21").Value = "1810"
$ws.Range("E22").Value = "1811"
$ws.Range("E23").Value = "1812"

# Remove the obsolete "2503" period row for DANIELA CAMILA FUENTES RAMOS (row 25),
# shifting the remaining rows up
$ws.Rows.Item(25).Delete()

# Row 24 (now the last data row of the table) needs the closing bottom border
# that used to belong to the removed row
$ws.Range("B24:J24").Borders.Item(9).LineStyle = 1
$ws.Range("B24:J24").Borders.Item(9).Weight = 2
$ws.Range("B24:J24").Borders.Item(9).ColorIndex = 1
